$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 538290.94
$ws.Range("I15").Value = 538290.94
$ws.Range("K15").Value = 1614872.82
$ws.Range("M15").Value = -1614703.82
$ws.Range("H18").Value = 1866.3334
$ws.Range("I18").Value = 1866.3334
$ws.Range("J18").Value = 0.0
$ws.Range("K18").Value = 1866.3334
$ws.Range("L18").Value = 0.0
$ws.Range("M18").ClearContents() | Out-Null
$ws.Range("N18").Value = -1582.3334
$ws.Range("H19").Value = 1148.9286
$ws.Range("I19").Value = 1224.0
$ws.Range("J19").Value = 1013.8
$ws.Range("K19").Value = 1224.0
$ws.Range("L19").Value = 1013.8
$ws.Range("M19").Value = -1049.0
$ws.Range("N19").Value = -1363.8
$ws.Range("H28").Value = 1722.25
$ws.Range("I28").Value = 1963.4286
$ws.Range("J28").Value = 34.0
$ws.Range("K28").Value = 1963.4286
$ws.Range("L28").Value = 34.0
$ws.Range("M28").Value = -1478.4286
$ws.Range("N28").Value = -1004.0
$ws.Range("H33").Value = 242.4
$ws.Range("I33").Value = 235.9
$ws.Range("J33").Value = 255.4
$ws.Range("K33").Value = 235.9
$ws.Range("L33").Value = 255.4
$ws.Range("M33").Value = -6.900000000000006
$ws.Range("N33").Value = -713.4
$ws.Range("H64").Value = 14289885.0
$ws.Range("I64").Value = 22225932.0
$ws.Range("J64").Value = 5000.0
$ws.Range("K64").Value = 22225932.0
$ws.Range("L64").Value = 5000.0
$ws.Range("M64").Value = -22225684.0
$ws.Range("N64").Value = -5496.0
$ws.Range("H67").Value = 14289885.0
$ws.Range("I67").Value = 22225932.0
$ws.Range("J67").Value = 5000.0
$ws.Range("K67").Value = 22225932.0
$ws.Range("L67").Value = 5000.0
$ws.Range("M67").Value = -22225074.0
$ws.Range("N67").Value = -6716.0
$ws.Range("H80").Value = 500.46155
$ws.Range("I80").Value = 436.83334
$ws.Range("J80").Value = 555.0
$ws.Range("K80").Value = 1310.50002
$ws.Range("L80").Value = 1665.0
$ws.Range("M80").Value = -312.5000199999999
$ws.Range("N80").Value = -3661.0
$ws.Range("H83").Value = 500.46155
$ws.Range("I83").Value = 436.83334
$ws.Range("J83").Value = 555.0
$ws.Range("K83").Value = 3931.50006
$ws.Range("L83").Value = 4995.0
$ws.Range("M83").Value = 1060.49994
$ws.Range("N83").Value = -14979.0
$ws.Range("H88").Value = 28389174.0
$ws.Range("I88").Value = 83334136.0
$ws.Range("J88").Value = 3969192.2
$ws.Range("K88").Value = 83334136.0
$ws.Range("L88").Value = 3969192.2
$ws.Range("M88").Value = -83333730.0
$ws.Range("N88").Value = -3970004.2
$ws.Range("H91").Value = 28389174.0
$ws.Range("I91").Value = 83334136.0
$ws.Range("J91").Value = 3969192.2
$ws.Range("K91").Value = 83334136.0
$ws.Range("L91").Value = 3969192.2
$ws.Range("M91").Value = -83332732.0
$ws.Range("N91").Value = -3972000.2
$ws.Range("H98").Value = 1814.7587
$ws.Range("I98").Value = 1699.12
$ws.Range("K98").Value = 1699.12
$ws.Range("M98").Value = -201.1199999999999
$ws.Range("H122").Value = 1814.7587
$ws.Range("I122").Value = 1699.12
$ws.Range("K122").Value = 5097.36
$ws.Range("M122").Value = -2647.36
$ws.Range("H138").Value = 2161.8096
$ws.Range("I138").Value = 1718.0
$ws.Range("J138").Value = 2883.0
$ws.Range("K138").Value = 5154.0
$ws.Range("L138").Value = 8649.0
$ws.Range("M138").Value = -14.0
$ws.Range("N138").Value = -18929.0
$ws.Range("H141").Value = 852.96875
$ws.Range("I141").Value = 852.96875
$ws.Range("J141").Value = 0.0
$ws.Range("K141").Value = 2558.90625
$ws.Range("L141").Value = 0.0
$ws.Range("M141").ClearContents() | Out-Null
$ws.Range("N141").Value = 2621.09375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4287.8965
$ws.Range("I32").Value = 4885.2173
$ws.Range("J32").Value = 1998.1666
$ws.Range("K32").Value = 4885.2173
$ws.Range("L32").Value = 1998.1666
$ws.Range("M32").Value = -4598.2173
$ws.Range("N32").Value = -2572.1666
$ws.Range("H61").Value = 58826796.0
$ws.Range("I61").Value = 71431680.0
$ws.Range("K61").Value = 71431680.0
$ws.Range("M61").Value = -71431468.0
$ws.Range("H62").Value = 46416.332
$ws.Range("J62").Value = 46416.332
$ws.Range("L62").Value = 46416.332
$ws.Range("N62").Value = -47664.332
$ws.Range("H63").Value = 3684.125
$ws.Range("I63").Value = 3684.125
$ws.Range("J63").Value = 0.0
$ws.Range("K63").Value = 3684.125
$ws.Range("L63").Value = 0.0
$ws.Range("M63").ClearContents() | Out-Null
$ws.Range("N63").Value = -2998.125
$ws.Range("H65").Value = 46416.332
$ws.Range("J65").Value = 46416.332
$ws.Range("L65").Value = 139248.996
$ws.Range("N65").Value = -145488.996
$ws.Range("H66").Value = 3684.125
$ws.Range("I66").Value = 3684.125
$ws.Range("J66").Value = 0.0
$ws.Range("K66").Value = 18420.625
$ws.Range("L66").Value = 0.0
$ws.Range("M66").ClearContents() | Out-Null
$ws.Range("N66").Value = -14988.625
$ws.Range("H97").Value = 799.2
$ws.Range("I97").Value = 799.2
$ws.Range("K97").Value = 799.2
$ws.Range("M97").Value = -303.2
$ws.Range("H132").Value = 2705280.2
$ws.Range("I132").Value = 2859500.2
$ws.Range("J132").Value = 6431.5
$ws.Range("K132").Value = 8578500.600000001
$ws.Range("L132").Value = 19294.5
$ws.Range("M132").Value = -8575970.600000001
$ws.Range("N132").Value = -24354.5
$ws.Range("H136").Value = 58826796.0
$ws.Range("I136").Value = 71431680.0
$ws.Range("K136").Value = 214295040.0
$ws.Range("M136").Value = -214292490.0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2008.9048
$ws.Range("I99").Value = 1848.875
$ws.Range("K99").Value = 1848.875
$ws.Range("M99").Value = -350.875
$ws.Range("H134").Value = 39617624.0
$ws.Range("I134").Value = 39617624.0
$ws.Range("K134").Value = 118852872.0
$ws.Range("M134").Value = -118850337.0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 12499.5
$ws.Range("I21").Value = 12499.5
$ws.Range("K21").Value = 12499.5
$ws.Range("M21").Value = -12264.5
$ws.Range("H33").Value = 2199.75
$ws.Range("I33").Value = 2199.75
$ws.Range("J33").Value = 0.0
$ws.Range("K33").Value = 2199.75
$ws.Range("L33").Value = 0.0
$ws.Range("M33").ClearContents() | Out-Null
$ws.Range("N33").Value = -1820.75
$ws.Range("H50").Value = 0.0
$ws.Range("J50").Value = 0.0
$ws.Range("L50").ClearContents() | Out-Null
$ws.Range("N50").Value = 0.0
$ws.Range("H51").Value = 19999.0
$ws.Range("J51").Value = 0.0
$ws.Range("L51").Value = 0.0
$ws.Range("N51").ClearContents() | Out-Null
$ws.Range("H59").Value = 68816.664
$ws.Range("I59").Value = 7633.6665
$ws.Range("J59").Value = 129999.664
$ws.Range("K59").Value = 7633.6665
$ws.Range("L59").Value = 129999.664
$ws.Range("M59").Value = -6488.6665
$ws.Range("N59").Value = -132289.664
$ws.Range("H60").Value = 27758.766
$ws.Range("J60").Value = 40909.0
$ws.Range("L60").Value = 40909.0
$ws.Range("N60").Value = -41931.0
$ws.Range("H61").Value = 19999.0
$ws.Range("J61").Value = 0.0
$ws.Range("L61").Value = 0.0
$ws.Range("N61").ClearContents() | Out-Null
$ws.Range("H99").Value = 3755.5
$ws.Range("I99").Value = 3755.5
$ws.Range("K99").Value = 3755.5
$ws.Range("M99").Value = -2257.5
$ws.Range("H126").Value = 3755.5
$ws.Range("I126").Value = 3755.5
$ws.Range("K126").Value = 11266.5
$ws.Range("M126").Value = -8796.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6720.0
$ws.Range("J80").Value = 3442.5
$ws.Range("L80").Value = 10327.5
$ws.Range("N80").Value = -12199.5
$ws.Range("H83").Value = 6720.0
$ws.Range("J83").Value = 3442.5
$ws.Range("L83").Value = 30982.5
$ws.Range("N83").Value = -40342.5
$ws.Range("H92").Value = 594.0
$ws.Range("I92").Value = 399.5
$ws.Range("J92").Value = 788.5
$ws.Range("K92").Value = 1198.5
$ws.Range("L92").Value = 2365.5
$ws.Range("M92").Value = 49.5
$ws.Range("N92").Value = -4861.5
$ws.Range("H123").Value = 5242.375
$ws.Range("J123").Value = 6684.8335
$ws.Range("L123").Value = 20054.5005
$ws.Range("N123").Value = -24954.5005
$ws.Range("H131").Value = 1720.8334
$ws.Range("I131").Value = 1088.3334
$ws.Range("K131").Value = 3265.0002
$ws.Range("M131").Value = 1774.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 99999.0
$ws.Range("J123").Value = 99999.0
$ws.Range("L123").Value = 99999.0
$ws.Range("N123").Value = -104899.0
$ws.Range("H132").Value = 9620102.0
$ws.Range("I132").Value = 9620102.0
$ws.Range("K132").Value = 28860306.0
$ws.Range("M132").Value = -28857776.0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 33333.0
$ws.Range("I26").Value = 33333.0
$ws.Range("K26").Value = 33333.0
$ws.Range("M26").Value = -33038.0
$ws.Range("H61").Value = 3058.4285
$ws.Range("J61").Value = 3003.0
$ws.Range("L61").Value = 3003.0
$ws.Range("N61").Value = -3407.0
$ws.Range("H113").Value = 3058.4285
$ws.Range("J113").Value = 3003.0
$ws.Range("L113").Value = 3003.0
$ws.Range("N113").Value = -7343.0
$ws.Range("H131").Value = 100000.0
$ws.Range("I131").Value = 0.0
$ws.Range("K131").Value = 0.0
$ws.Range("M131").ClearContents() | Out-Null
$ws.Range("H132").Value = 32010260.0
$ws.Range("I132").Value = 48012940.0
$ws.Range("J132").Value = 4899.4
$ws.Range("K132").Value = 144038820.0
$ws.Range("L132").Value = 14698.2
$ws.Range("M132").Value = -144036290.0
$ws.Range("N132").Value = -19758.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 13475.333
$ws.Range("I34").Value = 13475.333
$ws.Range("K34").Value = 13475.333
$ws.Range("M34").Value = -13272.333
$ws.Range("H132").Value = 19237406.0
$ws.Range("J132").Value = 100000.0
$ws.Range("L132").Value = 300000.0
$ws.Range("N132").Value = -305060.0
